$d = $word.ActiveDocument

# Locate the paragraph that contains the Savigny/Volksgeist answer.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Volksgeist*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph containing 'Volksgeist'"
}

# Range covering just the paragraph's content (exclude the trailing
# paragraph mark character) so the paragraph properties/identity are kept.
$pStart = $target.Range.Start
$pEnd = $target.Range.End - 1
$r = $d.Range($pStart, $pEnd)

# Replace the single run with three runs, wrapping "Volksgeist" in
# proofErr spellStart/spellEnd markers exactly as Word's proofer would
# when it flags an unrecognised word, splitting the sentence around it.
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r><w:t xml:space="preserve">A: Savigny argued that law grows with the people and should evolve naturally from customs and </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Volksgeist</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> (spirit of the people), not be imposed artificially.</w:t></w:r>' +
  '</w:p>'

$r.InsertXML($xmlFrag)
